$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the confidence-interval columns produced by the rerun.
$ws.Range("G1").Value = "ci.lower"
$ws.Range("H1").Value = "ci.upper"

# New ci.lower / ci.upper values for each row of the results table.
$ws.Range("G2").Value = -0.378027274878099
$ws.Range("H2").Value = 0.0954932243005093

$ws.Range("G3").Value = -0.117043562835675
$ws.Range("H3").Value = 0.0438877631572555

$ws.Range("G4").Value = -0.101732995390858
$ws.Range("H4").Value = 0.0555873114048591

$ws.Range("G5").Value = -0.343279409379245
$ws.Range("H5").Value = 0.128719299469298

$ws.Range("G6").Value = -0.307014589606692
$ws.Range("H6").Value = 0.167753987118282

$ws.Range("G7").Value = -0.348178542047175
$ws.Range("H7").Value = 0.130556324667412

$ws.Range("G8").Value = -0.302633012487934
$ws.Range("H8").Value = 0.165359876035549

$ws.Range("G9").Value = -0.108585868047373
$ws.Range("H9").Value = 0.0489351262151642

$ws.Range("G10").Value = -0.0783298249321993
$ws.Range("H10").Value = -0.0326826032378995

$ws.Range("G11").Value = -0.0361392207717058
$ws.Range("H11").Value = 0.00912910507928552

# Rows 12-14 (q_b11_b21, q_rxy1_rxy2, cross_over_point) have no CI values,
# matching their existing blank C:F cells.

$ws.Range("G15").Value = -0.0964906119045593
$ws.Range("H15").Value = 0.0638499857647706

$ws.Range("G16").Value = -0.28703808872716
$ws.Range("H16").Value = 0.189939492738466

$ws.Range("G17").Value = -0.295863328140885
$ws.Range("H17").Value = 0.194252479363243

$ws.Range("G18").Value = -0.203465990781716
$ws.Range("H18").Value = 0.111174622809718

$ws.Range("G19").Value = -0.605266024975867
$ws.Range("H19").Value = 0.330719752071099

$ws.Range("G20").Value = -0.614029179213385
$ws.Range("H20").Value = 0.335507974236564

$ws.Range("G21").Value = -0.00912910507928552
$ws.Range("H21").Value = 0.0361392207717058

$ws.Range("G22").Value = -0.0978702524303284
$ws.Range("H22").Value = 0.217171736094746

$ws.Range("G23").Value = -0.0271570551997267
$ws.Range("H23").Value = 0.107506136127106

$ws.Range("G24").Value = -0.291142212142182
$ws.Range("H24").Value = 0.64603756597433

$ws.Range("G25").Value = -0.0294852518213235
$ws.Range("H25").Value = 0.10478475924286

$ws.Range("G26").Value = -0.291698532135467
$ws.Range("H26").Value = 0.645519244533824
